# Auto-generated script to append new scrim rows to several sheets,
# matching the target diff (new registration rows logged on 2025-07-10).

$wb = $excel.ActiveWorkbook

# Colors (matches existing conditional team-color formatting already on the sheets):
#   columns A-C (brawlers 1-3): light blue fill
#   columns D-F (brawlers 4-6): light red/pink fill
#   column G (winner/team):      same fill as A-C or D-F depending on team, bold font
#   columns H-N (players/time):  no fill
$fillBlue = 16770508   # RGB(0xCC,0xE5,0xFF)
$fillRed  = 13421812   # RGB(0xF4,0xCC,0xCC)

function Set-ScrimRow {
    param($ws, $rowNum, $values)
    for ($c = 1; $c -le 14; $c++) {
        $cell = $ws.Cells.Item($rowNum, $c)
        $cell.Value = $values[$c - 1]
        $cell.Borders.LineStyle = 1
        if ($c -ge 1 -and $c -le 3) {
            $cell.Interior.Color = $fillBlue
        } elseif ($c -ge 4 -and $c -le 6) {
            $cell.Interior.Color = $fillRed
        } elseif ($c -eq 7) {
            $cell.Font.Bold = $true
            if ($values[6] -eq 'Equipo 1') {
                $cell.Interior.Color = $fillBlue
            } else {
                $cell.Interior.Color = $fillRed
            }
        }
    }
}

# --- Sneaky Fields (sheet1.xml) ---
$ws = $wb.Worksheets.Item('Sneaky Fields')
Set-ScrimRow $ws 76 @('GUS', 'SHADE', 'SPIKE', 'SANDY', 'DRACO', 'MAX', 'Equipo 2', 'SK|Yoshi825🇱🇻', 'SK|Joker', 'SK|OPE🏒', 'IC|Mebius', 'IC|RamaZR', 'IC|Nob', '20250710T134644.000Z')
Set-ScrimRow $ws 77 @('GUS', 'SHADE', 'SPIKE', 'SANDY', 'DRACO', 'MAX', 'Equipo 2', 'SK|Yoshi825🇱🇻', 'SK|Joker', 'SK|OPE🏒', 'IC|Mebius', 'IC|RamaZR', 'IC|Nob', '20250710T134307.000Z')
Set-ScrimRow $ws 78 @('GUS', 'SHADE', 'SPIKE', 'SANDY', 'DRACO', 'MAX', 'Equipo 1', 'SK|Yoshi825🇱🇻', 'SK|Joker', 'SK|OPE🏒', 'IC|Mebius', 'IC|RamaZR', 'IC|Nob', '20250710T134139.000Z')

# --- Crystal Arcade (sheet14.xml) ---
$ws = $wb.Worksheets.Item('Crystal Arcade')
Set-ScrimRow $ws 56 @('JANET', 'CHARLIE', 'LUMI', 'TARA', 'BONNIE', 'DARRYL', 'Equipo 2', 'RC|Battoman', 'RC|Shu', 'nyamura', 'ZETA|Sizuku', 'ZETA|Levi', 'ZETA|Sitetampo', '20250710T134710.000Z')
Set-ScrimRow $ws 57 @('JANET', 'CHARLIE', 'LUMI', 'TARA', 'BONNIE', 'DARRYL', 'Equipo 2', 'RC|Battoman', 'RC|Shu', 'nyamura', 'ZETA|Sizuku', 'ZETA|Levi', 'ZETA|Sitetampo', '20250710T134524.000Z')
Set-ScrimRow $ws 58 @('JANET', 'CHARLIE', 'LUMI', 'TARA', 'BONNIE', 'DARRYL', 'Equipo 1', 'RC|Battoman', 'RC|Shu', 'nyamura', 'ZETA|Sizuku', 'ZETA|Levi', 'ZETA|Sitetampo', '20250710T134255.000Z')

# --- Dry Season (sheet16.xml) ---
$ws = $wb.Worksheets.Item('Dry Season')
Set-ScrimRow $ws 62 @('DARRYL', 'LUMI', 'KIT', 'CARL', 'CHARLIE', 'GENE', 'Equipo 2', 'CR|Moya', 'CR|Milkreo', 'Tensai 천재', 'NAVI|Ryohei', 'NAVI|Achapi', 'NAVI|Kuru', '20250710T135544.000Z')
Set-ScrimRow $ws 63 @('DARRYL', 'LUMI', 'KIT', 'CARL', 'CHARLIE', 'GENE', 'Equipo 1', 'CR|Moya', 'CR|Milkreo', 'Tensai 천재', 'NAVI|Ryohei', 'NAVI|Achapi', 'NAVI|Kuru', '20250710T135324.000Z')
Set-ScrimRow $ws 64 @('DARRYL', 'LUMI', 'KIT', 'CARL', 'CHARLIE', 'GENE', 'Equipo 2', 'CR|Moya', 'CR|Milkreo', 'Tensai 천재', 'NAVI|Ryohei', 'NAVI|Achapi', 'NAVI|Kuru', '20250710T135123.000Z')
Set-ScrimRow $ws 65 @('DRACO', 'GENE', 'CHARLIE', 'ALLI', 'SQUEAK', 'BELLE', 'Equipo 1', 'CR|Moya', 'Tensai 천재', 'CR|Milkreo', 'NAVI|Ryohei', 'NAVI|Achapi', 'NAVI|Kuru', '20250710T134544.000Z')
Set-ScrimRow $ws 66 @('DRACO', 'GENE', 'CHARLIE', 'ALLI', 'SQUEAK', 'BELLE', 'Equipo 1', 'CR|Moya', 'Tensai 천재', 'CR|Milkreo', 'NAVI|Ryohei', 'NAVI|Achapi', 'NAVI|Kuru', '20250710T134336.000Z')
Set-ScrimRow $ws 67 @('GENE', 'CHARLIE', 'BROCK', 'MR. P', 'MANDY', 'JAE-YONG', 'Equipo 1', 'TTM|Angelboy', 'TTM|Maury', 'TTM|Maru', 'KDS|Decaii', 'KDS|Remica', 'KDS|Ćiro', '20250710T134305.000Z')
Set-ScrimRow $ws 68 @('GENE', 'CHARLIE', 'BROCK', 'MR. P', 'MANDY', 'JAE-YONG', 'Equipo 1', 'TTM|Angelboy', 'TTM|Maury', 'TTM|Maru', 'KDS|Decaii', 'KDS|Remica', 'KDS|Ćiro', '20250710T134045.000Z')

# --- Goldarm Gulch (sheet18.xml) ---
$ws = $wb.Worksheets.Item('Goldarm Gulch')
Set-ScrimRow $ws 66 @('GUS', 'OLLIE', 'BROCK', 'WILLOW', 'BELLE', 'ALLI', 'Equipo 1', 'nyamura', 'RC|Battoman', 'RC|Shu', 'ZETA|Levi', 'ZETA|Sizuku', 'ZETA|Sitetampo', '20250710T135648.000Z')
Set-ScrimRow $ws 67 @('GUS', 'OLLIE', 'BROCK', 'WILLOW', 'BELLE', 'ALLI', 'Equipo 1', 'nyamura', 'RC|Battoman', 'RC|Shu', 'ZETA|Levi', 'ZETA|Sizuku', 'ZETA|Sitetampo', '20250710T135347.000Z')
Set-ScrimRow $ws 68 @('LUMI', 'BROCK', 'KAZE', 'HANK', 'WILLOW', 'JAE-YONG', 'Equipo 1', 'CR|Moya', 'CR|Milkreo', 'Tensai 천재', 'NAVI|Ryohei', 'NAVI|Achapi', 'NAVI|Kuru', '20250710T140147.000Z')
Set-ScrimRow $ws 69 @('BELLE', 'OLLIE', 'BROCK', 'HANK', 'GUS', 'PIPER', 'Equipo 2', 'nyamura', 'RC|Battoman', 'RC|Shu', 'ZETA|Levi', 'ZETA|Sitetampo', 'ZETA|Sizuku', '20250710T140208.000Z')

# --- Open Business (sheet2.xml) ---
$ws = $wb.Worksheets.Item('Open Business')
Set-ScrimRow $ws 46 @('HANK', 'BEA', 'ALLI', 'LUMI', 'DRACO', 'SHADE', 'Equipo 2', 'GEN|Moding', 'GEN|cookie', 'GEN|BONOX2', 'RVL|Terry', 'Isee in Bush', 'RVL|Mameshi', '20250710T140027.000Z')
Set-ScrimRow $ws 47 @('ASH', 'SHADE', 'LOU', 'KAZE', 'DRACO', 'LUMI', 'Equipo 1', 'GEN|Moding', 'GEN|BONOX2', 'GEN|cookie', 'Isee in Bush', 'RVL|Terry', 'RVL|Mameshi', '20250710T135405.000Z')
Set-ScrimRow $ws 48 @('ASH', 'SHADE', 'LOU', 'KAZE', 'DRACO', 'LUMI', 'Equipo 1', 'GEN|Moding', 'GEN|BONOX2', 'GEN|cookie', 'Isee in Bush', 'RVL|Terry', 'RVL|Mameshi', '20250710T135219.000Z')
Set-ScrimRow $ws 49 @('ASH', 'SHADE', 'LOU', 'KAZE', 'DRACO', 'LUMI', 'Equipo 2', 'GEN|Moding', 'GEN|BONOX2', 'GEN|cookie', 'Isee in Bush', 'RVL|Terry', 'RVL|Mameshi', '20250710T135051.000Z')

# --- Triple Dribble (sheet7.xml) ---
$ws = $wb.Worksheets.Item('Triple Dribble')
Set-ScrimRow $ws 81 @('LUMI', 'WILLOW', 'ASH', 'KAZE', 'HANK', 'BERRY', 'Equipo 1', 'GEN|cookie', 'GEN|Moding', 'GEN|BONOX2', 'Isee in Bush', 'RVL|Terry', 'RVL|Mameshi', '20250710T134525.000Z')
Set-ScrimRow $ws 82 @('LUMI', 'WILLOW', 'ASH', 'KAZE', 'HANK', 'BERRY', 'Equipo 1', 'GEN|cookie', 'GEN|Moding', 'GEN|BONOX2', 'Isee in Bush', 'RVL|Terry', 'RVL|Mameshi', '20250710T134316.000Z')

# --- Belle's Rock (sheet8.xml) ---
$ws = $wb.Worksheets.Item('Belle''s Rock')
Set-ScrimRow $ws 81 @('TICK', 'JAE-YONG', 'HANK', 'DARRYL', 'BROCK', 'RICO', 'Equipo 1', 'NOVO|Marco', 'NOVO|Biso', 'NOVO|26is', 'TH|Code: LeNain', 'TH|Zhar', 'TH|IKAUSSA', '20250710T134120.000Z')

Write-Host "done"
